# Fruta / hortaliza, semanal
# Insert a new week's worth of Chirimoya price data (6 rows) for the
# "Vega Central Mapocho de Santiago" market, pushing the existing data
# down and expanding the used range from T110 to T116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows above the current row 95 (row formatting, e.g. the
# date style on column D, is inherited from the row being pushed down).
$ws.Rows("95:100").Insert()

function Set-Row {
    param($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q, $rr, $s, $t)
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j
    $ws.Range("K$r").Value = $k
    $ws.Range("L$r").Value = $l
    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p
    $ws.Range("Q$r").Value = $q
    $ws.Range("R$r").Value = $rr
    $ws.Range("S$r").Value = $s
    $ws.Range("T$r").Value = $t
}

Set-Row 95  9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Cuarta"                 220 1200  1200  1200  "`$/kilo (en caja de 15 kilos)" "Provincia de Limarí" 1200 1
Set-Row 96  9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Especial"               330 20000 20000 20000 "`$/bandeja 8 kilos"            "Provincia de Limarí" 2500 8
Set-Row 97  9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Extra (doble especial)" 300 24000 24000 24000 "`$/bandeja 8 kilos"            "Provincia de Limarí" 3000 8
Set-Row 98  9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Primera"                280 16000 16000 16000 "`$/bandeja 8 kilos"            "Provincia de Limarí" 2000 8
Set-Row 99  9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Segunda"                250 14400 14400 14400 "`$/bandeja 8 kilos"            "Provincia de Limarí" 1800 8
Set-Row 100 9 "Vega Central Mapocho de Santiago" "Metropolitana" 44504 13 "Fruta" 100107 "Otros" 100107002 "Chirimoya" "Cultivar IV Región" "Tercera"                300 1500  1500  1500  "`$/kilo (en caja de 15 kilos)" "Provincia de Limarí" 1500 1
